$wb = $excel.ActiveWorkbook

# 1) Update the "Last Updated" timestamp on the Metadata sheet.
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "29 Oct 2025, 09:24 PM"

# 2) Add the new "distance from Dma50" sheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "distance from Dma50"

# 3) Header row - reuse the header formatting already used by the other
#    sheets (bold, thin box border, centered) by copying it across.
$srcHeader = $wb.Worksheets.Item("1 Month Performance").Range("A1:C1")
$srcHeader.Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)

$ws.Range("A1").Value = "Icon"
$ws.Range("B1").Value = "Stock"
$ws.Range("C1").Value = "Distance From Sma50"

# 4) Data rows.
$data = @(
    ,@("📈", "NIFTYPSUBANK", 10.2033)
    ,@("📈", "NIFTYMETAL", 8.624700000000001)
    ,@("📈", "NIFTYOILANDGAS", 6.396)
    ,@("📈", "NIFTYCOMMODITIES", 5.7207)
    ,@("📈", "CNXINFRA", 5.6012)
    ,@("📈", "CNXREALTY", 5.4493)
    ,@("📈", "NIFTYPVTBANK", 5.0059)
    ,@("📈", "BANKNIFTY", 4.9192)
    ,@("📈", "NIFTYFINSERVICE", 3.9783)
    ,@("📈", "NIFTYMIDCAP50", 3.9228)
    ,@("📈", "NIFTY", 3.7191)
    ,@("📈", "CNXENERGY", 3.706)
    ,@("📈", "CNXMIDCAP", 3.6313)
    ,@("📈", "NIFTY200", 3.5915)
    ,@("📈", "NIFTY100", 3.5759)
    ,@("📈", "NIFTY500", 3.3379)
    ,@("📈", "CNXSMALLCAP", 2.9205)
    ,@("📈", "NIFTY50VALUE20", 2.893)
    ,@("📈", "NIFTYCPSE", 2.837)
    ,@("📈", "CNXNIFTYJUNIOR", 2.8291)
    ,@("📈", "NIFTYHEALTHCARE", 2.162)
    ,@("📈", "CNXIT", 2.0641)
    ,@("📈", "NIFTYCONSUMPTION", 2.0245)
    ,@("📈", "CNXPHARMA", 1.573)
    ,@("📈", "NIFTYAUTO", 1.5538)
    ,@("📈", "NIFTYGROWSECT15", 1.5325)
    ,@("📈", "NIFTYFMCG", 1.3194)
    ,@("📈", "NIFTYCONSURDURBL", 0.4031)
    ,@("📈", "NIFTYMEDIA", -1.9217)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

Write-Output ("Added sheet '" + $ws.Name + "' with " + $data.Count + " data rows.")
